# Presales executive-presentation.pptx trim + formatting fix
#
# 1. Remove all slides except the Title slide, the "Why This Solution?" slide,
#    the "Business Value - Financial Impact" slide, and the "Risk Mitigation"
#    slide (final deck goes from 17 slides down to 4).
# 2. Un-bold a handful of table cells on the three surviving table slides.

$p = $ppt.ActivePresentation

# --- 1. Trim the deck down to slides 1, 7, 8, 11 (original slide order) ----
# Delete from the highest index down to the lowest so indices of the slides
# we still need to delete don't shift underneath us.
$deletePositions = @(17, 16, 15, 14, 13, 12, 10, 9, 6, 5, 4, 3, 2)
foreach ($pos in $deletePositions) {
    $p.Slides.Item($pos).Delete()
}

# After the deletions above, the deck is:
#   1 -> Title slide (unchanged)
#   2 -> "Slide 7: Why This Solution?"
#   3 -> "Slide 8: Business Value - Financial Impact"
#   4 -> "Slide 11: Risk Mitigation"

function Unbold-Cell($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Font.Bold = $false
}

function Get-Table($slide) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.HasTable) {
            return $sh.Table
        }
    }
    return $null
}

# --- 2a. Slide 2 ("Why This Solution?") -------------------------------------
# Row 4 ([Current limitation 3] / [Our advantage 3]) loses its bold.
$tbl2 = Get-Table($p.Slides.Item(2))
Unbold-Cell $tbl2 4 1
Unbold-Cell $tbl2 4 2

# --- 2b. Slide 3 ("Business Value - Financial Impact") ---------------------
# Header row (Metric / Value) and the ROI row lose their bold.
$tbl3 = Get-Table($p.Slides.Item(3))
Unbold-Cell $tbl3 1 1
Unbold-Cell $tbl3 1 2
Unbold-Cell $tbl3 6 1
Unbold-Cell $tbl3 6 2

# --- 2c. Slide 4 ("Risk Mitigation") ----------------------------------------
# Header row (Risk / Mitigation Strategy / Success Probability) and the
# [Risk 3] row lose their bold.
$tbl4 = Get-Table($p.Slides.Item(4))
Unbold-Cell $tbl4 1 1
Unbold-Cell $tbl4 1 2
Unbold-Cell $tbl4 1 3
Unbold-Cell $tbl4 4 1
Unbold-Cell $tbl4 4 2
Unbold-Cell $tbl4 4 3
